$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Scanner" to "Session"
$ws.Name = "Session"

# Delete row 152 (A152=211926, D152=11:51:08) and shift cells up
$ws.Rows(152).Delete()
